# This script updates the single-column results table in the document.
# Rows are addressed by (1-based) row index within the one table in the
# document, so that cells sharing identical text values elsewhere in the
# table are not accidentally affected.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "825"

$t.Cell(6, 1).Range.Text  = "0.00066"
$t.Cell(7, 1).Range.Text  = "0.00020"
$t.Cell(8, 1).Range.Text  = "0.00006"
$t.Cell(9, 1).Range.Text  = "0.00025"
$t.Cell(10, 1).Range.Text = "0.00034"
$t.Cell(11, 1).Range.Text = "0.00048"
$t.Cell(12, 1).Range.Text = "0.16465"

$t.Cell(44, 1).Range.Text = "99.91"
$t.Cell(45, 1).Range.Text = "0.16"
$t.Cell(46, 1).Range.Text = "193"
